$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data values
$ws.Range("A2").Value = 5.6
$ws.Range("C2").Value = 9
$ws.Range("C3").Value = 8

# Change the formula in C4 from AVERAGE to SUM
$ws.Range("C4").Formula = "=SUM(C2:C3)"

# Add the new value in C6
$ws.Range("C6").Value = 9

# Update the selection to match the new active cell
[void]$ws.Range("I16").Select()
